# Update 'nombre_aides' (col C) and 'montant_total' (col D) figures for the
# 2020-07-28 Fonds de solidarite volet 2 data refresh (regional x categorie juridique).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "830"
$ws.Range("C4").Style = "Normal"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "2484355.40"
$ws.Range("D4").Style = "Normal"

$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "147"
$ws.Range("C20").Style = "Normal"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "370569.00"
$ws.Range("D20").Style = "Normal"

$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = "302"
$ws.Range("C21").Style = "Normal"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "960217.92"
$ws.Range("D21").Style = "Normal"

$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "147"
$ws.Range("C22").Style = "Normal"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "391137.26"
$ws.Range("D22").Style = "Normal"

$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "205"
$ws.Range("C28").Style = "Normal"

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "482926.00"
$ws.Range("D28").Style = "Normal"

$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "413"
$ws.Range("C30").Style = "Normal"

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1378959.16"
$ws.Range("D30").Style = "Normal"

$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = "321"
$ws.Range("C32").Style = "Normal"

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "950862.96"
$ws.Range("D32").Style = "Normal"

$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "59"
$ws.Range("C40").Style = "Normal"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "158871.00"
$ws.Range("D40").Style = "Normal"

$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "48"
$ws.Range("C41").Style = "Normal"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "226729.92"
$ws.Range("D41").Style = "Normal"

$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "75"
$ws.Range("C42").Style = "Normal"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "346547.99"
$ws.Range("D42").Style = "Normal"

$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "256"
$ws.Range("C44").Style = "Normal"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "712676.74"
$ws.Range("D44").Style = "Normal"

$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "495"
$ws.Range("C46").Style = "Normal"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1718187.88"
$ws.Range("D46").Style = "Normal"

$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "325"
$ws.Range("C47").Style = "Normal"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1013453.79"
$ws.Range("D47").Style = "Normal"

$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "24"
$ws.Range("C49").Style = "Normal"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "119240.09"
$ws.Range("D49").Style = "Normal"

$ws.Range("C70").NumberFormat = "@"
$ws.Range("C70").Value = "812"
$ws.Range("C70").Style = "Normal"

$ws.Range("D70").NumberFormat = "@"
$ws.Range("D70").Value = "2580905.34"
$ws.Range("D70").Style = "Normal"

$ws.Range("C71").NumberFormat = "@"
$ws.Range("C71").Value = "467"
$ws.Range("C71").Style = "Normal"

$ws.Range("D71").NumberFormat = "@"
$ws.Range("D71").Value = "1380386.03"
$ws.Range("D71").Style = "Normal"

$ws.Range("C73").NumberFormat = "@"
$ws.Range("C73").Value = "29"
$ws.Range("C73").Style = "Normal"

$ws.Range("D73").NumberFormat = "@"
$ws.Range("D73").Value = "100231.09"
$ws.Range("D73").Style = "Normal"

$ws.Range("C80").NumberFormat = "@"
$ws.Range("C80").Value = "185"
$ws.Range("C80").Style = "Normal"

$ws.Range("D80").NumberFormat = "@"
$ws.Range("D80").Value = "425993.00"
$ws.Range("D80").Style = "Normal"

$ws.Range("C82").NumberFormat = "@"
$ws.Range("C82").Value = "455"
$ws.Range("C82").Style = "Normal"

$ws.Range("D82").NumberFormat = "@"
$ws.Range("D82").Value = "1460362.50"
$ws.Range("D82").Style = "Normal"

$ws.Range("C83").NumberFormat = "@"
$ws.Range("C83").Value = "171"
$ws.Range("C83").Style = "Normal"

$ws.Range("D83").NumberFormat = "@"
$ws.Range("D83").Value = "481976.09"
$ws.Range("D83").Style = "Normal"

$ws.Range("C84").NumberFormat = "@"
$ws.Range("C84").Value = "15"
$ws.Range("C84").Style = "Normal"

$ws.Range("D84").NumberFormat = "@"
$ws.Range("D84").Value = "45500.00"
$ws.Range("D84").Style = "Normal"

$ws.Range("C85").NumberFormat = "@"
$ws.Range("C85").Value = "5"
$ws.Range("C85").Style = "Normal"

$ws.Range("D85").NumberFormat = "@"
$ws.Range("D85").Value = "18670.00"
$ws.Range("D85").Style = "Normal"

$ws.Range("C86").NumberFormat = "@"
$ws.Range("C86").Value = "451"
$ws.Range("C86").Style = "Normal"

$ws.Range("D86").NumberFormat = "@"
$ws.Range("D86").Value = "1035224.67"
$ws.Range("D86").Style = "Normal"

$ws.Range("C88").NumberFormat = "@"
$ws.Range("C88").Value = "964"
$ws.Range("C88").Style = "Normal"

$ws.Range("D88").NumberFormat = "@"
$ws.Range("D88").Value = "2982988.02"
$ws.Range("D88").Style = "Normal"

$ws.Range("C90").NumberFormat = "@"
$ws.Range("C90").Value = "885"
$ws.Range("C90").Style = "Normal"

$ws.Range("D90").NumberFormat = "@"
$ws.Range("D90").Value = "2451868.45"
$ws.Range("D90").Style = "Normal"
